# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 312
    4  = 10325
    6  = 943
    7  = 1281
    8  = 7045
    9  = 16
    10 = 437
    11 = 200
    13 = 3180
    15 = 312
    16 = 648
    17 = 123
    18 = 866
    20 = 63
    21 = 1624
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
